# Auto-generated script applying value updates described by the diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 851.9231
$ws.Range("I41").Value = 521.4286
$ws.Range("J41").Value = 1237.5
$ws.Range("K41").Value = 521.4286
$ws.Range("L41").Value = 1237.5
$ws.Range("M41").Value = -81.42859999999996
$ws.Range("N41").Value = -2117.5
$ws.Range("H43").Value = 1150.8889
$ws.Range("I43").Value = 1149.5
$ws.Range("J43").Value = 1151.2858
$ws.Range("K43").Value = 1149.5
$ws.Range("L43").Value = 1151.2858
$ws.Range("M43").Value = -1080.5
$ws.Range("N43").Value = -1289.2858
$ws.Range("H94").Value = 4445
$ws.Range("I94").Value = 3519.1667
$ws.Range("K94").Value = 3519.1667
$ws.Range("M94").Value = -3068.1667
$ws.Range("H97").Value = 20703.1
$ws.Range("J97").Value = 20703.1
$ws.Range("L97").Value = 62109.3
$ws.Range("N97").Value = -63101.3
$ws.Range("H135").Value = 572.1818
$ws.Range("I135").Value = 319
$ws.Range("J135").Value = 1247.3334
$ws.Range("K135").Value = 2871
$ws.Range("L135").Value = 11226.0006
$ws.Range("M135").Value = -336
$ws.Range("N135").Value = -16296.0006
$ws.Range("H137").Value = 2780731
$ws.Range("I137").Value = 3131
$ws.Range("J137").Value = 4632464.5
$ws.Range("K137").Value = 9393
$ws.Range("L137").Value = 13897393.5
$ws.Range("M137").Value = -6843
$ws.Range("N137").Value = -13902493.5
$ws.Range("H141").Value = 2396.875
$ws.Range("I141").Value = 2557.6924
$ws.Range("J141").Value = 1700
$ws.Range("K141").Value = 7673.0772
$ws.Range("L141").Value = 5100
$ws.Range("M141").Value = -2493.0772
$ws.Range("N141").Value = -15460
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3774.875
$ws.Range("I61").Value = 2566.5
$ws.Range("J61").Value = 7400
$ws.Range("K61").Value = 2566.5
$ws.Range("L61").Value = 7400
$ws.Range("M61").Value = -2354.5
$ws.Range("N61").Value = -7824
$ws.Range("H74").Value = 45414.52
$ws.Range("I74").Value = 56813.832
$ws.Range("J74").Value = 16102
$ws.Range("K74").Value = 56813.832
$ws.Range("L74").Value = 16102
$ws.Range("M74").Value = -55939.832
$ws.Range("N74").Value = -17850
$ws.Range("H77").Value = 45414.52
$ws.Range("I77").Value = 56813.832
$ws.Range("J77").Value = 16102
$ws.Range("K77").Value = 284069.16
$ws.Range("L77").Value = 80510
$ws.Range("M77").Value = -279701.16
$ws.Range("N77").Value = -89246
$ws.Range("H122").Value = 1843.5
$ws.Range("I122").Value = 1463.7646
$ws.Range("J122").Value = 2340.077
$ws.Range("K122").Value = 4391.293799999999
$ws.Range("L122").Value = 7020.231000000001
$ws.Range("M122").Value = -1941.293799999999
$ws.Range("N122").Value = -11920.231
$ws.Range("H136").Value = 3774.875
$ws.Range("I136").Value = 2566.5
$ws.Range("J136").Value = 7400
$ws.Range("K136").Value = 7699.5
$ws.Range("L136").Value = 22200
$ws.Range("M136").Value = -5149.5
$ws.Range("N136").Value = -27300
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5963.316
$ws.Range("I134").Value = 6084.28
$ws.Range("J134").Value = 5730.6924
$ws.Range("K134").Value = 18252.84
$ws.Range("L134").Value = 17192.0772
$ws.Range("M134").Value = -15717.84
$ws.Range("N134").Value = -22262.0772
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 98.15385000000001
$ws.Range("I7").Value = 73.333336
$ws.Range("J7").Value = 119.42857
$ws.Range("K7").Value = 73.333336
$ws.Range("L7").Value = 119.42857
$ws.Range("M7").Value = 39.666664
$ws.Range("N7").Value = -345.42857
$ws.Range("H22").Value = 208.75
$ws.Range("I22").Value = 155
$ws.Range("J22").Value = 262.5
$ws.Range("K22").Value = 155
$ws.Range("L22").Value = 262.5
$ws.Range("M22").Value = 195
$ws.Range("N22").Value = -962.5
$ws.Range("H31").Value = 2290.2942
$ws.Range("I31").Value = 1417.0834
$ws.Range("J31").Value = 4386
$ws.Range("K31").Value = 1417.0834
$ws.Range("L31").Value = 4386
$ws.Range("M31").Value = -1122.0834
$ws.Range("N31").Value = -4976
$ws.Range("H34").Value = 2290.2942
$ws.Range("I34").Value = 1417.0834
$ws.Range("J34").Value = 4386
$ws.Range("K34").Value = 1417.0834
$ws.Range("L34").Value = 4386
$ws.Range("M34").Value = -1215.0834
$ws.Range("N34").Value = -4790
$ws.Range("H132").Value = 2652.761
$ws.Range("I132").Value = 2338.516
$ws.Range("J132").Value = 3302.2
$ws.Range("K132").Value = 7015.548000000001
$ws.Range("L132").Value = 9906.599999999999
$ws.Range("M132").Value = -4485.548000000001
$ws.Range("N132").Value = -14966.6
$ws.Range("H134").Value = 2068.0908
$ws.Range("I134").Value = 1836.1904
$ws.Range("J134").Value = 2473.9167
$ws.Range("K134").Value = 5508.5712
$ws.Range("L134").Value = 7421.750100000001
$ws.Range("M134").Value = -2973.5712
$ws.Range("N134").Value = -12491.7501
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 847.7
$ws.Range("J131").Value = 894.24445
$ws.Range("L131").Value = 2682.73335
$ws.Range("N131").Value = -12762.73335
$ws.Range("H136").Value = 55559804
$ws.Range("I136").Value = 83334700
$ws.Range("K136").Value = 250004100
$ws.Range("M136").Value = -249999000
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 33703.71
$ws.Range("I122").Value = 43081.207
$ws.Range("J122").Value = 1552.2858
$ws.Range("K122").Value = 129243.621
$ws.Range("L122").Value = 4656.857400000001
$ws.Range("M122").Value = -126793.621
$ws.Range("N122").Value = -9556.857400000001
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2706.9583
$ws.Range("I122").Value = 2211.889
$ws.Range("J122").Value = 3004
$ws.Range("K122").Value = 6635.667
$ws.Range("L122").Value = 9012
$ws.Range("M122").Value = -4185.667
$ws.Range("N122").Value = -13912
$ws.Range("H132").Value = 3443.2693
$ws.Range("I132").Value = 2392.6667
$ws.Range("J132").Value = 5807.125
$ws.Range("K132").Value = 7178.000100000001
$ws.Range("L132").Value = 17421.375
$ws.Range("M132").Value = -4648.000100000001
$ws.Range("N132").Value = -22481.375
$ws.Range("H136").Value = 1862.8125
$ws.Range("I136").Value = 1862.8125
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 5588.4375
$ws.Range("L136").Value = 0
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -3038.4375
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 44593.5
$ws.Range("I63").Value = 40226
$ws.Range("J63").Value = 46049.332
$ws.Range("K63").Value = 40226
$ws.Range("L63").Value = 46049.332
$ws.Range("M63").Value = -39602
$ws.Range("N63").Value = -47297.332
$ws.Range("H66").Value = 44593.5
$ws.Range("I66").Value = 40226
$ws.Range("J66").Value = 46049.332
$ws.Range("K66").Value = 120678
$ws.Range("L66").Value = 138147.996
$ws.Range("M66").Value = -117558
$ws.Range("N66").Value = -144387.996
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").ClearContents()
$ws.Range("N98").Value = 0
$ws.Range("H105").Value = 29333.334
$ws.Range("J105").Value = 29333.334
$ws.Range("L105").Value = 29333.334
$ws.Range("N105").Value = -36321.334
